# Adds the "Architecture Topic" subsection detail columns (F/G) to the
# Cross Reference Matrix sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content Creator block (rows 8-11): 1.3.1 / Creator Features
foreach ($r in 8..11) {
    $ws.Cells.Item($r, 6).Value = "1.3.1"
    $ws.Cells.Item($r, 7).Value = "Creator Features"
}

# Viewer block (rows 2-6): 1.3.2 / Viewer Features
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 6).Value = "1.3.2"
    $ws.Cells.Item($r, 7).Value = "Viewer Features"
}

# Advertisers block (rows 13-15): 1.3.3 / Advertiser Features
foreach ($r in 13..15) {
    $ws.Cells.Item($r, 6).Value = "1.3.3"
    $ws.Cells.Item($r, 7).Value = "Advertiser Features"
}

# NFR block (rows 17-21): 1.2 / Container Diagrams
foreach ($r in 17..21) {
    $ws.Cells.Item($r, 6).Value = 1.2
    $ws.Cells.Item($r, 7).Value = "Container Diagrams"
}

$ws.Range("F22").Select()
